# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 826
$ws1.Range("F6").Value = 282
$ws1.Range("F7").Value = 6731
$ws1.Range("F8").Value = 56
$ws1.Range("F15").Value = 20
$ws1.Range("F16").Value = 225
$ws1.Range("F17").Value = 576
$ws1.Range("F18").Value = 64

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 826
$ws4.Range("F6").Value = 282
$ws4.Range("F7").Value = 6731
$ws4.Range("F8").Value = 56
$ws4.Range("F15").Value = 20
$ws4.Range("F16").Value = 225
$ws4.Range("F17").Value = 576
$ws4.Range("F18").Value = 64
$ws4.Range("F19").Value = 7
